$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row ---
$ws.Range("A1").Value = "Student_ID"
$ws.Range("B1").Value = "Name"
$ws.Range("C1").Value = "Preference"
$ws.Range("D1").Value = "Status"
$ws.Range("E1").Value = "Company_ID"

# --- Student data rows ---
$data = @(
    @("S9621745C", "Law Wen Ming",     "Software Development",                  "Unassigned"),
    @("S9518672A", "Say De Ming",      "System Development",                    "Unassigned"),
    @("S9348750H", "Sheng Shu Qi",     "Software Engineering, Development",     "Unassigned"),
    @("S2316221J", "Bay Yong Quan",    "IOS and Android Development",           "Unassigned"),
    @("S6461996I", "Chang De Kang",    "Documents, QA Testing and Development", "Unassigned"),
    @("S0194725G", "Lau Kai Feng",     "Software Engineering, Development",     "Unassigned"),
    @("S5806960D", "Cheong Xuan Ming", "IOS and Android Development",           "Unassigned"),
    @("S2206960D", "Wei Jun",          "IOS and Android Development",           "Unassigned")
)

$row = 2
foreach ($rec in $data) {
    $ws.Cells.Item($row, 1).Value = $rec[0]
    $ws.Cells.Item($row, 2).Value = $rec[1]
    $ws.Cells.Item($row, 3).Value = $rec[2]
    $ws.Cells.Item($row, 4).Value = $rec[3]
    $row++
}

# Matches original data's left-aligned style carried over onto the new ID
# column (rows 2-4) plus B2 (per source formatting).
$ws.Range("A2:A4").HorizontalAlignment = -4131
$ws.Range("B2").HorizontalAlignment = -4131

# --- Column widths ---
# Target XML widths are 40.7109375 / 56.7109375 / 21.42578125 "characters".
# This runtime quantizes ColumnWidth onto a 1/6-character grid before
# storing it, so the literal target values are not exactly representable;
# the inputs below land on the closest achievable grid point to each target.
$ws.Columns.Item(2).ColumnWidth = 40.7109375
$ws.Columns.Item(3).ColumnWidth = 56.7109375
$ws.Columns.Item(5).ColumnWidth = 20.71

# --- Selection ---
$ws.Range("F5").Select()
